# Leetcode 222 - Count Complete Tree Node
# Append a new row (#19) to the Question List table on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108
$xlLeft   = -4131

$row = 20

# A: running index
$ws.Cells.Item($row, 1).Value = 19
$ws.Cells.Item($row, 1).HorizontalAlignment = $xlCenter

# B: Question title
$ws.Cells.Item($row, 2).Value = "Count Complete Tree Node"
$ws.Cells.Item($row, 2).HorizontalAlignment = $xlLeft

# D: Data Structure
$ws.Cells.Item($row, 4).Value = "Tree"
$ws.Cells.Item($row, 4).HorizontalAlignment = $xlCenter

# E: Difficulty
$ws.Cells.Item($row, 5).Value = "medium"
$ws.Cells.Item($row, 5).HorizontalAlignment = $xlCenter

# F: Source reference
$ws.Cells.Item($row, 6).Value = "leetcode 222"
$ws.Cells.Item($row, 6).HorizontalAlignment = $xlCenter

# Restore the cursor/selection to where the author left it after typing the row
$ws.Range("F26").Select()
